# Burndown chart and sprint log update
# Fill in the Day 5 (column I) actuals for the remaining sprint-log rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = 0.8
$ws.Range("I8").Value = 0.8
$ws.Range("I9").Value = 1
$ws.Range("I10").Value = 1

# Move the active selection to I6 (matches the author's cursor position
# after finishing the Day 5 entries).
$ws.Range("I6").Select()
